$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — match the existing header
# formatting (bold font, thin border, centered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2 (plain numeric values, default style)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
